$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New route strings (column A) for the 4 newly added endpoints
$ws.Range("A8").Value  = "Route::get('machines/{productionlineID}', FormController::class.'@getMachines');"
$ws.Range("A9").Value  = "Route::get('speedLosses/{PO}/{productionLine}', FormController::class.'@get_speedLosses');"
$ws.Range("A10").Value = "Route::get('pos/{shift}/{site}', FormController::class.'@getPOsFromShift');"
$ws.Range("A11").Value = "Route::get('events/{PO}/{productionLine}', FormController::class.'@getEvents');"

# Matching API-address formulas (column B), mirroring the existing rows' style
$ws.Range("B8").Formula  = '=_xlfn.CONCAT($B$2,"machine/{producitonLineId}")'
$ws.Range("B9").Formula  = '=_xlfn.CONCAT($B$2,"speedloss/{po}/{productionLine}")'
$ws.Range("B10").Formula = '=_xlfn.CONCAT($B$2,"pos/{shift}/{site}")'
$ws.Range("B11").Formula = '=_xlfn.CONCAT($B$2,"events/{po}/{productionline}")'

$excel.Calculate()

# Widen the two columns to fit the new, longer strings
$ws.Columns.Item(1).ColumnWidth = 82.16666666666667
$ws.Columns.Item(2).ColumnWidth = 46.83

# Move the active selection the same way the author's session ended up
$ws.Range("B23").Select()
